$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are stored as text, matching the
# original "inline string" representation used throughout column D.
$textRows = @(4,5,6,7,8,9,10,11,12,13,14,15,17,18,19,20,21,22,23,25,26,27,28,29,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51)
foreach ($r in $textRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "24.672.44"
$ws.Range("E2").Value = "  -4.53%  "

$ws.Range("D3").Value = "1.654.68"
$ws.Range("E3").Value = "  -5.23%  "

$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +1.39%  "

$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D5").Value = "1.005"
$ws.Range("E5").Value = "  +1.14%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "308.05"
$ws.Range("E6").Value = "  -2.50%  "

$ws.Range("D7").Value = "0.3642"
$ws.Range("E7").Value = "  -4.88%  "

$ws.Range("D8").Value = "0.3314"
$ws.Range("E8").Value = "  -8.94%  "

$ws.Range("D9").Value = "46.98"
$ws.Range("E9").Value = "  -7.19%  "

$ws.Range("D10").Value = "1.146"
$ws.Range("E10").Value = "  -6.60%  "

$ws.Range("D11").Value = "0.07179"
$ws.Range("E11").Value = "  -6.86%  "

$ws.Range("D12").Value = "1.009"
$ws.Range("E12").Value = "  +1.40%  "

$ws.Range("D13").Value = "6.086"
$ws.Range("E13").Value = "  -5.76%  "

$ws.Range("D14").Value = "20.06"
$ws.Range("E14").Value = "  -7.96%  "

$ws.Range("D15").Value = "6.720"
$ws.Range("E15").Value = "  -4.96%  "

$ws.Range("D16").Value = "1.651.09"
$ws.Range("E16").Value = "  -5.26%  "

$ws.Range("D17").Value = "0.00001077"
$ws.Range("E17").Value = "  -6.97%  "

$ws.Range("D18").Value = "1.005"
$ws.Range("E18").Value = "  +1.22%  "

$ws.Range("D19").Value = "0.06606"
$ws.Range("E19").Value = "  -3.26%  "

$ws.Range("D20").Value = "80.64"
$ws.Range("E20").Value = "  -7.46%  "

$ws.Range("D21").Value = "16.48"
$ws.Range("E21").Value = "  -6.13%  "

$ws.Range("D22").Value = "6.048"
$ws.Range("E22").Value = "  -6.51%  "

$ws.Range("D23").Value = "12.32"
$ws.Range("E23").Value = "  -3.51%  "

$ws.Range("D24").Value = "24.707.64"
$ws.Range("E24").Value = "  -4.10%  "

$ws.Range("D25").Value = "2.404"
$ws.Range("E25").Value = "  -1.14%  "

$ws.Range("D26").Value = "2.615"
$ws.Range("E26").Value = "  -11.06%  "

$ws.Range("D27").Value = "148.43"
$ws.Range("E27").Value = "  -3.84%  "

$ws.Range("D28").Value = "19.54"
$ws.Range("E28").Value = "  -5.58%  "

$ws.Range("D29").Value = "129.18"
$ws.Range("E29").Value = "  -3.83%  "

$ws.Range("D30").Value = "1.850.11"
$ws.Range("E30").Value = "  -4.31%  "

$ws.Range("D31").Value = "1.184"
$ws.Range("E31").Value = "  -1.16%  "

$ws.Range("D32").Value = "4.145"
$ws.Range("E32").Value = "  -4.56%  "

$ws.Range("D33").Value = "6.271"
$ws.Range("E33").Value = "  -11.00%  "

$ws.Range("D34").Value = "1.738"

$ws.Range("D35").Value = "0.08498"
$ws.Range("E35").Value = "  -2.30%  "

$ws.Range("D36").Value = "13.04"
$ws.Range("E36").Value = "  -8.96%  "

$ws.Range("D37").Value = "5.295"
$ws.Range("E37").Value = "  -6.21%  "

$ws.Range("D38").Value = "0.06294"
$ws.Range("E38").Value = "  -6.35%  "

$ws.Range("D39").Value = "0.02290"
$ws.Range("E39").Value = "  -6.95%  "

$ws.Range("D40").Value = "8.536"
$ws.Range("E40").Value = "  -8.28%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.2117"
$ws.Range("E41").Value = "  -4.72%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "1.230"
$ws.Range("E42").Value = "  -5.33%  "

$ws.Range("D43").Value = "0.6140"
$ws.Range("E43").Value = "  -6.19%  "

$ws.Range("D44").Value = "1.006"
$ws.Range("E44").Value = "  +1.29%  "

$ws.Range("D45").Value = "13.34"
$ws.Range("E45").Value = "  -4.17%  "

$ws.Range("D46").Value = "3.768"
$ws.Range("E46").Value = "  -3.40%  "

$ws.Range("D47").Value = "0.5838"
$ws.Range("E47").Value = "  -8.08%  "

$ws.Range("D48").Value = "2.003"
$ws.Range("E48").Value = "  -8.03%  "

$ws.Range("D49").Value = "123.88"
$ws.Range("E49").Value = "  -5.86%  "

$ws.Range("D50").Value = "0.07083"
$ws.Range("E50").Value = "  -5.40%  "

$ws.Range("D51").Value = "75.70"
$ws.Range("E51").Value = "  -4.67%  "
